$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.425.62"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "2.370.42"
$ws.Range("E3").Value = "  +5.40%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "232.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.95"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.26%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.457"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0951"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "2.721.66"
$ws.Range("E13").Value = "  +5.41%  "
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.843"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "2.370.71"
$ws.Range("D19").Value = "43.392.04"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "0.0₃0985"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +17.63%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.129"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.38%  "
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("E34").Value = "  +4.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0692"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.29%  "
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("E40").Value = "  -2.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.93%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  +3.69%  "
$ws.Range("E44").Value = "  +8.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0953"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").Value = "1.447.01"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.594.27"
$ws.Range("E50").Value = "  +5.61%  "
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000205"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.51%  "
